# Weekly update: a new price record (week of 2022-08-26) is inserted into the
# daily consolidated sheet for Acelga / Agrícola del Norte S.A. de Arica.
# This shifts every existing record from row 34 onward down by one row and
# grows the sheet dimension to A1:R71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 34, pushing rows 34..70 down to 35..71.
$ws.Rows("34").Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(34, 1).Value  = 1
$ws.Cells.Item(34, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value  = 44799
$ws.Cells.Item(34, 5).Value  = 15
$ws.Cells.Item(34, 6).Value  = 100112009
$ws.Cells.Item(34, 7).Value  = "Acelga"
$ws.Cells.Item(34, 8).Value  = "Sin especificar"
$ws.Cells.Item(34, 9).Value  = "Primera"
$ws.Cells.Item(34, 10).Value = 250
$ws.Cells.Item(34, 11).Value = 1000
$ws.Cells.Item(34, 12).Value = 1200
$ws.Cells.Item(34, 13).Value = 1100
$ws.Cells.Item(34, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 367
$ws.Cells.Item(34, 17).Value = 3
$ws.Cells.Item(34, 18).Value = "Hortaliza"
